$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "1.0000", "77.40").
# Force the range to Text format first so Excel keeps them exactly as
# typed instead of normalising them into numbers (stripping zeros, etc).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 5).Value = "  +0.40%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.745.21"
$ws.Cells.Item(3, 5).Value = "  +0.17%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.06%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "235.38"
$ws.Cells.Item(5, 5).Value = "  -0.24%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.08%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.5082"
$ws.Cells.Item(7, 5).Value = "  +3.26%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "40.38"
$ws.Cells.Item(8, 5).Value = "  -2.52%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.2655"
$ws.Cells.Item(9, 5).Value = "  +4.30%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "0.06147"
$ws.Cells.Item(10, 5).Value = "  +2.09%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "1.753.77"
$ws.Cells.Item(11, 5).Value = "  +0.67%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "0.06934"
$ws.Cells.Item(12, 5).Value = "  +1.56%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "15.30"
$ws.Cells.Item(13, 5).Value = "  +3.19%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "0.6204"
$ws.Cells.Item(14, 5).Value = "  +9.05%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "4.462"
$ws.Cells.Item(15, 5).Value = "  +0.53%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "77.40"
$ws.Cells.Item(16, 5).Value = "  +1.27%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "1.001"
$ws.Cells.Item(17, 5).Value = "  -0.05%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -0.04%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "25.802.51"
$ws.Cells.Item(19, 5).Value = "  +0.33%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "11.54"
$ws.Cells.Item(20, 5).Value = "  +2.15%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "0.000006585"
$ws.Cells.Item(21, 5).Value = "  +0.52%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "1.971.41"
$ws.Cells.Item(22, 5).Value = "  +0.31%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +0.75%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "8.233"
$ws.Cells.Item(24, 5).Value = "  +3.76%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "5.126"
$ws.Cells.Item(25, 5).Value = "  +1.52%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "136.03"
$ws.Cells.Item(26, 5).Value = "  -0.69%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "1.468"
$ws.Cells.Item(27, 5).Value = "  -0.60%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "15.00"
$ws.Cells.Item(28, 5).Value = "  +2.39%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -2.59%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "102.31"
$ws.Cells.Item(30, 5).Value = "  +0.57%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "0.08175"
$ws.Cells.Item(31, 5).Value = "  +2.60%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "3.660"
$ws.Cells.Item(32, 5).Value = "  -2.27%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "3.372"
$ws.Cells.Item(33, 5).Value = "  -0.62%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -0.30%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "Frax"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(35, 4).Value = "1.0000"
$ws.Cells.Item(35, 5).Value = "  -0.08%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "HuobiToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(36, 4).Value = "2.647"
$ws.Cells.Item(36, 5).Value = "  +1.28%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "ARBITRUM"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(37, 4).Value = "0.9910"
$ws.Cells.Item(37, 5).Value = "  +0.99%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "ImmutableX"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(38, 4).Value = "0.5951"
$ws.Cells.Item(38, 5).Value = "  -0.30%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "MXToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(39, 4).Value = "2.620"
$ws.Cells.Item(39, 5).Value = "  -1.48%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(40, 4).Value = "0.01552"
$ws.Cells.Item(40, 5).Value = "  +2.55%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "1.001"
$ws.Cells.Item(41, 5).Value = "  -0.10%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "RenderToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(42, 4).Value = "1.908"
$ws.Cells.Item(42, 5).Value = "  -0.27%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Quant"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(43, 4).Value = "101.48"
$ws.Cells.Item(43, 5).Value = "  -0.13%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "TheSandbox"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(44, 4).Value = "0.3809"
$ws.Cells.Item(44, 5).Value = "  +1.75%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "TrustWalletToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(45, 4).Value = "0.7449"
$ws.Cells.Item(45, 5).Value = "  -1.03%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "FraxShare"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(46, 4).Value = "4.870"
$ws.Cells.Item(46, 5).Value = "  -5.56%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Cronos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(47, 4).Value = "0.05485"
$ws.Cells.Item(47, 5).Value = "  +4.89%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "Algorand"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(48, 4).Value = "0.1091"
$ws.Cells.Item(48, 5).Value = "  +2.38%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Aptos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(49, 4).Value = "5.914"
$ws.Cells.Item(49, 5).Value = "  +2.09%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "Elrond"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(50, 4).Value = "29.96"
$ws.Cells.Item(50, 5).Value = "  -0.33%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Aave"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(51, 4).Value = "52.31"
$ws.Cells.Item(51, 5).Value = "  +0.46%  "
